$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)   # "总计"
$q3src = $wb.Worksheets.Item(2)   # existing "2022-Q3" sheet holding fund-level data

# --- Duplicate the existing "2022-Q3" sheet. The duplicate lands right after the
# --- original (position 3) and keeps the original's full content/format untouched,
# --- including the tabSelected flag on its sheetView.
$q3src.Copy($null, $q3src)

$q4 = $wb.Worksheets.Item(2)   # original sheet -> becomes "2022-Q4" (refreshed data)
$q3 = $wb.Worksheets.Item(3)   # duplicate      -> stays "2022-Q3" (untouched data)

$q4.Name = "2022-Q4"
$q3.Name = "2022-Q3"

# --- Update the "2022-Q4" sheet with the new quarter's figures -----------------

# Re-style the header row + first data-row label cell to match the bordered
# "总计" header style instead of the plain one the old sheet used.
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").PasteSpecial(-4122)

# The four text-looking numeric columns must stay text cells (not be coerced to
# numbers). Format as text, assign, then restore the plain (unstyled) format
# from a sibling cell so no stray style sticks to these cells.
$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "0.67"
$q4.Range("E2").Value = "91.81"
$q4.Range("F2").Value = "2.97"
$q4.Range("G2").Value = "0.0199"
$q4.Range("C2").Copy()
$q4.Range("D2:G2").PasteSpecial(-4122)

# Rank column is numeric.
$q4.Range("H2").Value = 7

# Match the page margins used elsewhere in the workbook (0.75in/1in/0.5in).
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# --- Update the "总计" summary sheet -------------------------------------------

# Row 2 now reports on 2022-Q4 instead of 2022-Q3.
$total.Range("B2").Value = "2022-Q4"

# Add a new row 3 carrying the prior 2022-Q3 summary, re-using A2's bordered style.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.02
